$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.5
$ws.Range("G2").Value = 3.85
$ws.Range("H2").Value = 2.2
$ws.Range("I2").Value = 2.28
$ws.Range("L2").Value = 1.45
$ws.Range("O2").Value = 1.38
$ws.Range("P2").Value = 1.81
$ws.Range("Q2").Value = 2.12
$ws.Range("R2").Value = 1.3
$ws.Range("T2").Value = 1.81
$ws.Range("V2").Value = 1.78
$ws.Range("W2").Value = 1.35
$ws.Range("X2").Value = 14
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 14.5
$ws.Range("AB2").Value = 13.5
$ws.Range("AE2").Value = 27
$ws.Range("AF2").Value = 25
$ws.Range("AK2").Value = 55
$ws.Range("AN2").Value = 55

# Row 3
$ws.Range("F3").Value = 2.02
$ws.Range("G3").Value = 2.08
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 4.4
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 3.75
$ws.Range("N3").Value = 3.8
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.93
$ws.Range("Q3").Value = 1.97
$ws.Range("R3").Value = 1.36
$ws.Range("S3").Value = 3.5
$ws.Range("T3").Value = 1.05
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.29
$ws.Range("W3").Value = 1.93
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 15.5
$ws.Range("Z3").Value = 32
$ws.Range("AA3").Value = 90
$ws.Range("AB3").Value = 9.4
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 17.5
$ws.Range("AE3").Value = 55
$ws.Range("AF3").Value = 12.5
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 65
$ws.Range("AJ3").Value = 24
$ws.Range("AL3").Value = 38
$ws.Range("AM3").Value = 110
$ws.Range("AN3").Value = 16
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 2.54
$ws.Range("G4").Value = 2.76
$ws.Range("I4").Value = 3.35
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 3.45
$ws.Range("L4").Value = 1.48
$ws.Range("M4").Value = 1.1
$ws.Range("P4").Value = 1.75
$ws.Range("Q4").Value = 2.16
$ws.Range("R4").Value = 1.28
$ws.Range("S4").Value = 3.95
$ws.Range("U4").Value = 1.99
$ws.Range("W4").Value = 1.57
$ws.Range("AF4").Value = 21
$ws.Range("AG4").Value = 12.5
$ws.Range("AK4").Value = 50
$ws.Range("AN4").Value = 32
$ws.Range("AO4").Value = 46

# Row 5
$ws.Range("F5").Value = 1.4
$ws.Range("G5").Value = 1.44
$ws.Range("H5").Value = 9
$ws.Range("I5").Value = 11
$ws.Range("J5").Value = 5.1
$ws.Range("K5").Value = 5.8
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 4.6
$ws.Range("O5").Value = 1.24
$ws.Range("P5").Value = 2.28
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 1.48
$ws.Range("S5").Value = 2.82
$ws.Range("T5").Value = 1.94
$ws.Range("U5").Value = 1.89
$ws.Range("V5").Value = 1.11
$ws.Range("W5").Value = 3.25
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 990
$ws.Range("AF5").Value = 11
$ws.Range("AG5").Value = 10
$ws.Range("AH5").Value = 990
$ws.Range("AJ5").Value = 14
$ws.Range("AK5").Value = 17.5
$ws.Range("AN5").Value = 6.2

# Row 6
$ws.Range("F6").Value = 3.2
$ws.Range("H6").Value = 2.4
$ws.Range("I6").Value = 2.44
$ws.Range("J6").Value = 3.6
$ws.Range("L6").Value = 1.41
$ws.Range("N6").Value = 3.75
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 1.96
$ws.Range("Q6").Value = 1.94
$ws.Range("S6").Value = 3.5
$ws.Range("U6").Value = 2.24
$ws.Range("V6").Value = 1.69
$ws.Range("W6").Value = 1.44
$ws.Range("X6").Value = 14
$ws.Range("AC6").Value = 9.199999999999999
$ws.Range("AI6").Value = 38
$ws.Range("AJ6").Value = 55
$ws.Range("AO6").Value = 21

# Row 7
$ws.Range("I7").Value = 13
$ws.Range("K7").Value = 6.6
$ws.Range("O7").Value = 1.2
$ws.Range("Q7").Value = 1.63
$ws.Range("S7").Value = 2.6
$ws.Range("T7").Value = 2.04
$ws.Range("U7").Value = 1.84
$ws.Range("W7").Value = 3.85
$ws.Range("Y7").Value = 44
$ws.Range("AA7").Value = 430
$ws.Range("AB7").Value = 9.4
$ws.Range("AE7").Value = 180
$ws.Range("AF7").Value = 8.800000000000001
$ws.Range("AG7").Value = 10.5
$ws.Range("AJ7").Value = 10.5
$ws.Range("AK7").Value = 13.5
$ws.Range("AL7").Value = 38
$ws.Range("AN7").Value = 4.9

# Row 8
$ws.Range("F8").Value = 1.78
$ws.Range("G8").Value = 1.84
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 3.6
$ws.Range("O8").Value = 1.64
$ws.Range("Q8").Value = 2.86
$ws.Range("R8").Value = 1.17
$ws.Range("S8").Value = 6.2
$ws.Range("T8").Value = 2.44
$ws.Range("AB8").Value = 5.9

# Row 9
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 2.3
$ws.Range("I9").Value = 4.1
$ws.Range("J9").Value = 3.35
$ws.Range("P9").Value = 2
$ws.Range("S9").Value = 3
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.33
$ws.Range("W9").Value = 1.78
